# Update notebook 3 (test_db_excel_w_ecoinvent) to ecoinvent 3.11:
#  - the biosphere exchange "database" cells (F4/F8) move from the old
#    "biosphere3" database name to "ecoinvent-3.11-biosphere"
#  - the diesel technosphere exchange (row 7) moves from
#    "ecoinvent 3.9 conseq" / its old activity code to
#    "ecoinvent-3.11-consequential" / the new activity code, and its Notes
#    cell is updated to the new (CH) activity name
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "ecoinvent-3.11-biosphere"
$ws.Range("F8").Value = "ecoinvent-3.11-biosphere"

$ws.Range("G7").Value = "5e00f9695a7ce345a4a17c517fd1ea62"
$ws.Range("F7").Value = "ecoinvent-3.11-consequential"

$ws.Range("O7").Value = "'diesel production, low-sulfur' (kilogram, CH, None)"

$ws.Range("A7").Select()
